# Update the "Data" sheet of MetaData RunMode1.xlsx:
# Replace the Cocci batch "20220919-Cocci-10427Updt" result rows (2-13) with the
# new Salm batch "20220928-Salm-14911Updt" (cartridge CartridgeSalm4911,
# result IDs A1013601-A1013612, lanes shifted to start at lane 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(
  @{Row=2;  ResultId="A1013606"; Lane="6"},
  @{Row=3;  ResultId="A1013607"; Lane="7"},
  @{Row=4;  ResultId="A1013608"; Lane="8"},
  @{Row=5;  ResultId="A1013609"; Lane="9"},
  @{Row=6;  ResultId="A1013610"; Lane="10"},
  @{Row=7;  ResultId="A1013611"; Lane="11"},
  @{Row=8;  ResultId="A1013612"; Lane="12"},
  @{Row=9;  ResultId="A1013601"; Lane="1"},
  @{Row=10; ResultId="A1013602"; Lane="2"},
  @{Row=11; ResultId="A1013603"; Lane="3"},
  @{Row=12; ResultId="A1013604"; Lane="4"},
  @{Row=13; ResultId="A1013605"; Lane="5"}
)

$labSampleId = "20220928-Salm-14911Updt"
$cartridgeId = "CartridgeSalm4911"

foreach ($item in $rows) {
  $r = $item.Row
  $ws.Range("A$r").Value = $item.ResultId
  $ws.Range("E$r").Value = $labSampleId
  $ws.Range("Q$r").Value = $item.Lane
  $ws.Range("T$r").Value = $cartridgeId
}
